$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new header cell in AC1 ("No Of Months"), extending the table by one
# column. The cell picks up the default (unstyled) format, matching the
# author's edit which left the new header without the bold-ish style used
# by R1:AB1.
$ws.Range("AC1").Value = "No Of Months"

# Give the new column a sensible, explicit width (13 characters) instead of
# the auto bestFit widths used by the existing columns.
$ws.Columns.Item(29).ColumnWidth = 12.1667

# Move the selection over to the newly added column / area, mirroring the
# reviewer scrolling the sheet to the right to check the new column after
# adding it.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 28
$ws.Range("AJ12").Select()
